{"js": "// Locate the paragraph that ends the \"01/03/2022 ...\" entry (the last\n// dated entry, referencing \"iterazione 2 / iterazione 3 / stesura\n// documentazione definitiva;\") and append three new dated entries after\n// it, matching the same paragraph formatting (hanging indent + it-IT\n// language) as the surrounding list items.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"iterazione 2 / iterazione 3 / stesura documentazione definitiva;\";\nlet anchor = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(marker) !== -1) {\n    anchor = p;\n  }\n}\nif (!anchor) {\n  throw new Error(\"Could not find the '01/03/2022' planning paragraph to anchor the new entries on.\");\n}\n\n// Each entry is written as two runs (date prefix, then the rest of the\n// sentence) to mirror how the original list items are authored.\nconst newEntries = [\n  [\"02/03/2022 - \", \"iterazione 3 / stesura documentazione definitiva;\"],\n  [\"03/03/2022 \", \"- iterazione 3 / stesura documentazione definitiva;\"],\n  [\"04/03/2022 \", \"- iterazione 3 / stesura documentazione definitiva;\"],\n];\n\nfor (const [first, second] of newEntries) {\n  anchor = anchor.insertParagraph(\"\", \"After\");\n  anchor.insertText(first, \"Start\");\n  anchor.insertText(second, \"End\");\n}\n\nawait context.sync();\n", "ps1": "# Find the paragraph containing the \"01/03/2022\" planning entry (the last\n# dated bullet, ending in \"iterazione 2 / iterazione 3 / stesura\n# documentazione definitiva;\") and add three new dated entries right\n# after it, matching the hanging-indent / it-IT formatting used by the\n# rest of the list.\n\n$d = $word.ActiveDocument\n\n$marker = \"iterazione 2 / iterazione 3 / stesura documentazione definitiva;\"\n\n$anchor = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*$marker*\") {\n        $anchor = $p\n    }\n}\nif ($anchor -eq $null) {\n    throw \"Could not find the '01/03/2022' planning paragraph to anchor the new entries on.\"\n}\n\n$newEntries = @(\n    , @(\"02/03/2022 - \", \"iterazione 3 / stesura documentazione definitiva;\")\n    , @(\"03/03/2022 \", \"- iterazione 3 / stesura documentazione definitiva;\")\n    , @(\"04/03/2022 \", \"- iterazione 3 / stesura documentazione definitiva;\")\n)\n\nforeach ($entry in $newEntries) {\n    $anchor.Range.InsertParagraphAfter()\n    $anchor = $anchor.Next()\n\n    $r = $anchor.Range\n    $r.InsertAfter($entry[0])\n\n    $r2 = $anchor.Range\n    $r2.MoveEnd(1, -1) | Out-Null\n    $r2.Collapse(0)\n    $r2.InsertAfter($entry[1])\n}\n"}
